# Auto-generated edit script: refresh cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.380.71"
$ws.Range("E2").Value = "  -4.42%  "
$ws.Range("D3").Value = "2.362.84"
$ws.Range("E3").Value = "  -6.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.21"
$ws.Range("E5").Value = "  -4.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.25"
$ws.Range("E6").Value = "  -3.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "2.378.96"
$ws.Range("E9").Value = "  -5.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0956"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.74"
$ws.Range("E12").Value = "  -9.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.315"
$ws.Range("E13").Value = "  -5.04%  "
$ws.Range("D14").Value = "2.789.44"
$ws.Range("E14").Value = "  -5.08%  "
$ws.Range("D15").Value = "56.316.32"
$ws.Range("E15").Value = "  -4.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.49"
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "2.378.99"
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.28"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.02"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.06"
$ws.Range("E21").Value = "  -4.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.92"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  -4.86%  "
$ws.Range("D27").Value = "2.471.43"
$ws.Range("E27").Value = "  -5.50%  "
$ws.Range("E28").Value = "  -5.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.50"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0712"
$ws.Range("E33").Value = "  -6.39%  "
$ws.Range("E34").Value = "  -7.17%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.68"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.19"
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.70"
$ws.Range("E39").Value = "  -6.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.65"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("E41").Value = "  -5.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.784"
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.89"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.83"
$ws.Range("E45").Value = "  -6.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "255.01"
$ws.Range("E46").Value = "  -7.99%  "
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0901"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0487"
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("E50").Value = "  -5.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.64"
$ws.Range("E51").Value = "  -6.26%  "

Write-Host "Updated cryptos list with GitHub Actions"
